$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.83 = 19071.5 pesos`n✅ 19071.5 pesos = 4.8 = 955.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet: rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 207
$wsTasas.Range("O10").Value = 3947.8
$wsTasas.Range("N12").Value = 3976
$wsTasas.Range("O12").Value = 199.24
